$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New release notes text for [1.12]
$version = "[1.12]"
$details = "[printing envelops]`n- add Spire to CDMailer`n- handle case when Print is clicked without setting the output folder`n- handle envelops as size A5`n- fix a bug in filename generation for both envelops and letters`n"
$releaseDate = 43281

# Row 12 becomes the new data row (style mirrors the row above it: s=3 / s=4 / s=5)
$ws.Range("A12").HorizontalAlignment = -4131
$ws.Range("A12").VerticalAlignment = -4160
$ws.Range("A12").Value = $version

$ws.Range("B12").HorizontalAlignment = -4131
$ws.Range("B12").VerticalAlignment = -4160
$ws.Range("B12").WrapText = $true
$ws.Range("B12").Value = $details

$ws.Range("C12").HorizontalAlignment = -4131
$ws.Range("C12").VerticalAlignment = -4160
$ws.Range("C12").NumberFormat = "d-mmm-yy"
$ws.Range("C12").Value = $releaseDate

$ws.Rows.Item(12).RowHeight = 90

# Row 13 is the new trailing blank row (style mirrors the old blank row 12: s=3/s=3/s=3)
$ws.Range("A13").HorizontalAlignment = -4131
$ws.Range("A13").VerticalAlignment = -4160
$ws.Range("B13").HorizontalAlignment = -4131
$ws.Range("B13").VerticalAlignment = -4160
$ws.Range("C13").HorizontalAlignment = -4131
$ws.Range("C13").VerticalAlignment = -4160

# Grow the table (and its autofilter) to include the new row
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C13"))
